$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G14").Value = 0.58
$ws.Range("G16").Value = 0.864
$ws.Range("G23").Value = 0.8120000000000001
$ws.Range("G26").Value = 0.58
$ws.Range("G27").Value = 0.7
$ws.Range("G34").Value = 0.7
$ws.Range("G42").Value = 0.58
$ws.Range("G44").Value = 0.58
$ws.Range("G47").Value = 0.9320000000000001
